# New crime data collected - update weekly CompStat figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: volume/issue number and the week-covering date range.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/16/2024  Through  9/22/2024"

# ---------------------------------------------------------------------
# Helper: convert a numeric "count" cell into the blank-style text cell
# used when a category has zero/suppressed activity for the period. The
# source workbook renders that with the shared string "0" (not the
# number 0) using the same font/format as the neighbouring text cells.
# Copying from a cell that is already in that state reproduces the
# exact style instead of minting a brand-new one.
# ---------------------------------------------------------------------
$zeroSrc = $ws.Range("D14")     # already s="14" t="s" v="0"
$suppressedSrc = $ws.Range("E14")  # already s="14" t="s" v="***.*"

$zeroSrc.Copy($ws.Range("C14"))
$zeroSrc.Copy($ws.Range("C15"))
$zeroSrc.Copy($ws.Range("D15"))
$suppressedSrc.Copy($ws.Range("E15"))
$zeroSrc.Copy($ws.Range("C22"))
$zeroSrc.Copy($ws.Range("C27"))
$zeroSrc.Copy($ws.Range("C28"))
$zeroSrc.Copy($ws.Range("C29"))
$zeroSrc.Copy($ws.Range("C30"))

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
$ws.Range("F15").Value = 5
$ws.Range("H15").Value = 400
$ws.Range("L15").Value = -17.948717948717
$ws.Range("M15").Value = 3.225806451612
$ws.Range("N15").Value = -41.818181818181

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = 12
$ws.Range("E16").Value = 8.333333333333
$ws.Range("F16").Value = 46
$ws.Range("G16").Value = 38
$ws.Range("H16").Value = 21.052631578947
$ws.Range("I16").Value = 383
$ws.Range("J16").Value = 344
$ws.Range("K16").Value = 11.337209302325
$ws.Range("L16").Value = 18.575851393188
$ws.Range("M16").Value = 28.093645484949
$ws.Range("N16").Value = -62.340216322517

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 25
$ws.Range("D17").Value = 17
$ws.Range("E17").Value = 47.058823529411
$ws.Range("F17").Value = 66
$ws.Range("G17").Value = 77
$ws.Range("H17").Value = -14.285714285714
$ws.Range("I17").Value = 581
$ws.Range("J17").Value = 619
$ws.Range("K17").Value = -6.138933764135
$ws.Range("L17").Value = 7.195571955719
$ws.Range("M17").Value = 88.025889967637
$ws.Range("N17").Value = -10.752688172043

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -5.882352941176
$ws.Range("I18").Value = 199
$ws.Range("J18").Value = 203
$ws.Range("K18").Value = -1.970443349753
$ws.Range("L18").Value = -5.238095238095
$ws.Range("M18").Value = -17.768595041322
$ws.Range("N18").Value = -85.537790697674

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 23
$ws.Range("E19").Value = 43.75
$ws.Range("F19").Value = 73
$ws.Range("G19").Value = 63
$ws.Range("H19").Value = 15.873015873015
$ws.Range("I19").Value = 661
$ws.Range("J19").Value = 564
$ws.Range("K19").Value = 17.198581560283
$ws.Range("L19").Value = 19.314079422382
$ws.Range("M19").Value = 187.391304347826
$ws.Range("N19").Value = 64.427860696517

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
$ws.Range("D20").Value = 17
$ws.Range("E20").Value = -47.058823529411
$ws.Range("F20").Value = 50
$ws.Range("G20").Value = 47
$ws.Range("H20").Value = 6.382978723404
$ws.Range("I20").Value = 400
$ws.Range("J20").Value = 470
$ws.Range("K20").Value = -14.893617021276
$ws.Range("L20").Value = 23.839009287925
$ws.Range("M20").Value = 70.212765957446
$ws.Range("N20").Value = -66.244725738396

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 72
$ws.Range("D21").Value = 64
$ws.Range("E21").Value = 12.5
$ws.Range("F21").Value = 257
$ws.Range("G21").Value = 243
$ws.Range("H21").Value = 5.761316872427
$ws.Range("I21").Value = 2261
$ws.Range("J21").Value = 2241
$ws.Range("K21").Value = 0.892458723784
$ws.Range("L21").Value = 12.937062937062
$ws.Range("M21").Value = 65.641025641025
$ws.Range("N21").Value = -51.975361087510

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 19
$ws.Range("K22").Value = 15.789473684210
$ws.Range("L22").Value = -29.032258064516

# ---------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = -30
$ws.Range("I23").Value = 80
$ws.Range("J23").Value = 80
$ws.Range("L23").Value = -5.882352941176
$ws.Range("M23").Value = 48.148148148148

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = 3.448275862068
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 98
$ws.Range("H24").Value = -2.040816326530
$ws.Range("I24").Value = 945
$ws.Range("J24").Value = 983
$ws.Range("K24").Value = -3.865717192268
$ws.Range("L24").Value = -14.711191335740
$ws.Range("M24").Value = 71.818181818181

# ---------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 166.666666666667
$ws.Range("F25").Value = 36
$ws.Range("G25").Value = 25
$ws.Range("H25").Value = 44
$ws.Range("I25").Value = 286
$ws.Range("J25").Value = 274
$ws.Range("K25").Value = 4.379562043795
$ws.Range("L25").Value = -14.880952380952

# ---------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 19
$ws.Range("D26").Value = 27
$ws.Range("E26").Value = -29.629629629629
$ws.Range("G26").Value = 72
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 824
$ws.Range("J26").Value = 696
$ws.Range("K26").Value = 18.390804597701
$ws.Range("L26").Value = 25.801526717557
$ws.Range("M26").Value = 15.568022440392

# ---------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 50
$ws.Range("J27").Value = 47
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = -27.692307692307

# ---------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = -22.222222222222
$ws.Range("J28").Value = 58
$ws.Range("K28").Value = 12.068965517241
$ws.Range("L28").Value = 22.641509433962

# ---------------------------------------------------------------------
# Row 29 - Shooting Vic.
# ---------------------------------------------------------------------
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = -100
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = -66.666666666666
$ws.Range("J29").Value = 34
$ws.Range("K29").Value = -29.411764705882
$ws.Range("N29").Value = -77.358490566037

# ---------------------------------------------------------------------
# Row 30 - Shooting Inc.
# ---------------------------------------------------------------------
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = -100
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = -66.666666666666
$ws.Range("J30").Value = 29
$ws.Range("K30").Value = -44.827586206896
$ws.Range("N30").Value = -84.158415841584
